# [TID] Add a NOp to remove the single instruction loop, fix counters,
# empty queue even when waiting for memory.
#
# Fills in the previously-blank "Loose Pipeline" / "Loose Pipeline Pre HDD"
# rows (16-21, which used to be #DIV/0! because C/D were blank) and appends
# two new six-row blocks: "Rogue Instruction Removed" (rows 23-25) and
# "Rogue Instruction Removed HDD" (rows 26-28), each with its own
# AVERAGE() roll-up in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Existing "Loose Pipeline" block (rows 16-18) and "Loose Pipeline Pre HDD"
# block (rows 19-21) — these rows already existed but only had the A/B
# shared formula; fill in C:H with real data so B stops being #DIV/0!.
# ---------------------------------------------------------------------

$looseData = @(
    @(2768085, 926439, 622526, 164, 291748, 137),
    @(2470527, 725114, 722729, 188, 296580, 159),
    @(3275411, 1073800, 823286, 220, 303451, 181),
    @(2340704, 926307, 418593,  80,  68854,  88),
    @(2042036, 724977, 518584, 100,  72811, 110),
    @(2843368, 1073647, 618575, 120,  76768, 132)
)

for ($i = 0; $i -lt $looseData.Length; $i++) {
    $row = 16 + $i
    $vals = $looseData[$i]
    $ws.Cells.Item($row, 3).Value = $vals[0]   # C
    $ws.Cells.Item($row, 4).Value = $vals[1]   # D
    $ws.Cells.Item($row, 5).Value = $vals[2]   # E
    $ws.Cells.Item($row, 6).Value = $vals[3]   # F
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
    $ws.Cells.Item($row, 8).Value = $vals[5]   # H
}

# Extend the B6:B21 shared formula down through B28 (skipping the blank
# spacer rows 8, 15 and 22 that separate each block, which never had a
# formula), and add the I21 roll-up average for the Loose Pipeline blocks.
$cpiRows = @(6, 7, 9, 10, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 23, 24, 25, 26, 27, 28)
foreach ($r in $cpiRows) {
    $ws.Range("B$r").Formula = "=`$C$r/`$D$r"
}
$ws.Range("I21").Formula = "=AVERAGE(B16:B21)"

# ---------------------------------------------------------------------
# New "Rogue Instruction Removed" block (row 23 header, rows 23-25) and
# "Rogue Instruction Removed HDD" block (row 26 header, rows 26-28).
# ---------------------------------------------------------------------

$ws.Range("A23").Value = "Rogue Instruction Removed"
$ws.Range("A26").Value = "Rogue Instruction Removed HDD"

$rogueData = @(
    @(2670527,  8403, 621514, 172, 2027926, 147),
    @(3970749, 13810, 721224, 188, 3214266, 171),
    @(3875633, 12933, 821352, 220, 3021147, 197),
    @(2240704,  8260, 417227,  80, 1802363,  98),
    @(3542258, 13673, 517079, 100, 2990497, 122),
    @(3443590, 12780, 616641, 120, 2794464, 148)
)

for ($i = 0; $i -lt $rogueData.Length; $i++) {
    $row = 23 + $i
    $vals = $rogueData[$i]
    $ws.Cells.Item($row, 3).Value = $vals[0]   # C
    $ws.Cells.Item($row, 4).Value = $vals[1]   # D
    $ws.Cells.Item($row, 5).Value = $vals[2]   # E
    $ws.Cells.Item($row, 6).Value = $vals[3]   # F
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
    $ws.Cells.Item($row, 8).Value = $vals[5]   # H
}

$ws.Range("I28").Formula = "=AVERAGE(B23:B28)"

# ---------------------------------------------------------------------
# Selection matches the new end-of-data focus cell.
# ---------------------------------------------------------------------

$ws.Range("H27").Select() | Out-Null
